$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# H2: 30-Nov-2021 -> 8-Jan-2022 (keep existing d-mmm format)
$ws.Range("H2").Value = 44569

# H5: 11-Dec-2021 -> 27-Jan-2022 (keep existing d-mmm format)
$ws.Range("H5").Value = 44588

# H6: "20 hojas" -> "20 HOJAS"
$ws.Range("H6").Value = "20 HOJAS"

# H7: 8-Jan-2022 -> 12-Feb-2022 (keep existing d-mmm format)
$ws.Range("H7").Value = 44604

# H8: keep text "20 HOJAS", but now gets the d-mmm-yy number format (style 6)
$ws.Range("H8").NumberFormat = "d-mmm-yy"
$ws.Range("H8").Value = "20 HOJAS"

# H9: new cell with text "8 MARZO.,"
$ws.Range("H9").Value = "8 MARZO.,"

# H10: was a date (27-Jan-2022) -> becomes text "20 HOJAS", keeping the d-mmm format
$ws.Range("H10").Value = "20 HOJAS"

# H11: was text "20 HOJAS" (General format) -> becomes date 25-Mar-2022 with the d-mmm format
$ws.Range("H11").NumberFormat = "d-mmm"
$ws.Range("H11").Value = 44645

# H12: was a date (12-Feb-2022) -> becomes text "20 hojas", keeping the d-mmm format
$ws.Range("H12").Value = "20 hojas"

# H13, H14, H15: fully cleared (cell removed, including style)
$ws.Range("H13").Clear()
$ws.Range("H14").Clear()
$ws.Range("H15").Clear()

# Update the active selection to H13 (matches the saved view state)
$ws.Range("H13").Select()
